# "No-stock checks on database" - append five new "5 - 10 <Shop>" sheets
# (PullAndBear, Mango, Zara, Stradivarius, Bershka), each a 2-row stock
# check log (header + single "no changes yet" row), mirroring the layout
# already used by the other per-run sheets in this workbook.
$wb = $excel.ActiveWorkbook

function Add-NoStockSheet {
    param(
        [string]$Name,
        [string]$Hora
    )
    $lastIndex = $wb.Worksheets.Count
    $src = $wb.Worksheets.Item("1 - 10 Bershka")
    $src.Copy($null, $wb.Worksheets.Item($lastIndex))
    $new = $wb.Worksheets.Item($lastIndex + 1)
    $new.Name = $Name
    $new.Rows("3:13").Delete()
    $new.Range("A2").Value = $Hora
    $new.Range("B2").Value = $false
    $new.Range("C2").Value = 0
    $new.Range("D2").Value = 0
}

Add-NoStockSheet "5 - 10 PullAndBear" "23:6"
Add-NoStockSheet "5 - 10 Mango" "23:6"
Add-NoStockSheet "5 - 10 Zara" "23:7"
Add-NoStockSheet "5 - 10 Stradivarius" "23:7"
Add-NoStockSheet "5 - 10 Bershka" "23:7"
